# Generate Report for Handoff
# The b527d46b-... source file has moved from "In Translation" to
# "Ready for handoff" for both the zh-cn and de-de locales. Update the
# per-locale detail sheets with the new status + handoff timestamp, and
# roll the same status/timestamp up into the Overview summary sheet.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 3 is the b527d46b-... entry ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-17 20:13:57"

# --- de-de sheet: row 3 is the b527d46b-... entry ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-17 20:14:00"

# --- Overview sheet: row 3 rolls up the b527d46b-... entry ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-14-17 20:14:00"
